# "change for multi order" -- add Instructions / Mode of Payment / Total / Split
# columns (N:Q) with data for the two existing shipment rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("N1").Value = "Instructions"
$ws.Range("O1").Value = "Mode of Payment"
$ws.Range("P1").Value = "Total"
$ws.Range("Q1").Value = "Split"

# Row 2 (first shipment)
$ws.Range("N2").Value = "try"
$ws.Range("O2").Value = "paid"
$ws.Range("P2").Value = 23
$ws.Range("Q2").Value = 55

# Row 3 (second shipment)
$ws.Range("N3").Value = "again"
$ws.Range("O3").Value = "paid"
$ws.Range("P3").Value = 43
$ws.Range("Q3").Value = 89

# Keep the view roughly consistent with where the new data now lives.
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("P15").Select()
